# "moved local logs to repo" - fill in this week's task summary data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Cumulative total label now carries the computed total inline
$ws.Range("A14").Value = "Cumulative Total: 80"

# Row 3: Project analysis task
$ws.Range("A3").Value = "Project analysis"

# Row 4: Project design task
$ws.Range("A4").Value = "Project design"

$ws.Range("B3").Value = "analyse requirments and develop feature set for first build"

$ws.Range("B4").Value = "discuss various aspects of development/design"

# Header: author name and week number
$ws.Range("C1").Value = "Jesse Hare"
$ws.Range("E1").Value = 4

$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 10

$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 10

# Leave the selection where the author last left it
$ws.Range("D6").Select()
